$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the bold/bordered header style (currently on B1) onto A1, since in the
# new layout A1 becomes the first header cell and needs that same style.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column A on the data rows (2-3) currently carries the header style (from the old
# "4"/"18" values) - that formatting should not survive the shift, so clear it.
$ws.Range("A2:A3").Clear()

# Column F is no longer used after the shift.
$ws.Range("F1:F3").Clear()

# Row 1 (headers) - each column moved one letter to the left.
$ws.Range("A1").Value = "QS_Astral_exact5"
$ws.Range("B1").Value = "FNRATE_ASTRAL"
$ws.Range("C1").Value = "TAXON"
$ws.Range("D1").Value = "MODEL_CONDITION"
$ws.Range("E1").Value = "GENE"

# Row 2
$ws.Range("A2").Value = 1520
$ws.Range("B2").Value = 0.125
$ws.Range("C2").Value = "11-texon"
$ws.Range("D2").Value = "estimated_5genes_weakILS"
$ws.Range("E2").Value = 4

# Row 3
$ws.Range("A3").Value = 1520
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "11-texon"
$ws.Range("D3").Value = "estimated_5genes_weakILS"
$ws.Range("E3").Value = 18
